$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update swapped numeric values in row 18
$ws.Range("D18").Value = 100
$ws.Range("I18").Value = 89

# Update swapped text values in row 20
$ws.Range("I20").Value = "ind5"
$ws.Range("K20").Value = "ind3"

# Update sheet view (scroll position and active selection)
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 10
$ws.Range("I21").Select()
